$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: write column A (the period labels) in the precise order
# required so the six brand-new shared strings land at shared-string
# table indexes 21-26 (existing labels keep their original slots no
# matter when they are rewritten).
$ws.Cells.Item(3, 1).Value = "September 28, 1932 - February 27, 1933"
$ws.Cells.Item(4, 1).Value = "February 21, 1934 - March 14, 1935"
$ws.Cells.Item(5, 1).Value = "March 10, 1937 - March 31, 1938"
$ws.Cells.Item(6, 1).Value = "September 9, 1939 - June 5, 1940"
$ws.Cells.Item(7, 1).Value = "November 12, 1940 - April 28, 1942"
$ws.Cells.Item(2, 1).Value = "September 16, 1929 - June 1, 1932"
$ws.Cells.Item(8, 1).Value = "May 29, 1946 - October 9, 1946"
$ws.Cells.Item(9, 1).Value = "December 12, 1961 - June 26, 1962"
$ws.Cells.Item(10, 1).Value = "February 9, 1966 - October 7, 1966"
$ws.Cells.Item(11, 1).Value = "November 29, 1968 - May 26, 1970"
$ws.Cells.Item(12, 1).Value = "January 11, 1973 - October 3, 1974"
$ws.Cells.Item(13, 1).Value = "September 21, 1976 - February 28, 1978"
$ws.Cells.Item(14, 1).Value = "November 28, 1980 - August 12, 1982"
$ws.Cells.Item(15, 1).Value = "August 25, 1987 - December 4, 1987"
$ws.Cells.Item(16, 1).Value = "July 16, 1990 - October 11, 1990"
$ws.Cells.Item(17, 1).Value = "July 17, 1998 - August 31, 1998"
$ws.Cells.Item(18, 1).Value = "March 24, 2000 - October 9, 2002"
$ws.Cells.Item(19, 1).Value = "October 9, 2007 - March 9, 2009"
$ws.Cells.Item(20, 1).Value = "May 2, 2011 - October 4, 2011"
$ws.Cells.Item(21, 1).Value = "September 20, 2018 - December 24, 2018"
$ws.Cells.Item(22, 1).Value = "February 19, 2020 - March 23, 2020"
$ws.Cells.Item(23, 1).Value = "January 3, 2022 - October 12, 2022"

# --- Step 2: write Peak/Trough/Decline/Duration for every data row (2-23).
$rows = @(
    @{ Row = 2; B = 31.71; BFmt = $true; C = 4.4000000000000004; CFmt = $true; D = 0.8619; E = 679 },
    @{ Row = 3; B = 8.36; BFmt = $true; C = 5.53; CFmt = $true; D = 0.40600000000000003; E = 102 },
    @{ Row = 4; B = 11.43; BFmt = $true; C = 8.06; CFmt = $true; D = 0.2984; E = 264 },
    @{ Row = 5; B = 18.670000000000002; BFmt = $true; C = 8.5; CFmt = $true; D = 0.54469999999999996; E = 268 },
    @{ Row = 6; B = 13.17; BFmt = $true; C = 9.09; CFmt = $true; D = 0.30980000000000002; E = 183 },
    @{ Row = 7; B = 11.36; BFmt = $true; C = 7.47; CFmt = $true; D = 0.34239999999999998; E = 150 },
    @{ Row = 8; B = 19.25; BFmt = $true; C = 15.75; CFmt = $true; D = 0.182; E = 133 },
    @{ Row = 9; B = 72.64; BFmt = $true; C = 52.32; CFmt = $true; D = 0.28000000000000003; E = 196 },
    @{ Row = 10; B = 94.059997999999993; BFmt = $true; C = 69.290001000000004; CFmt = $true; D = 0.22; E = 240 },
    @{ Row = 11; B = 108.370003; BFmt = $true; C = 69.290001000000004; CFmt = $true; D = 0.36099999999999999; E = 543 },
    @{ Row = 12; B = 120.239998; BFmt = $true; C = 62.279998999999997; CFmt = $true; D = 0.48; E = 630 },
    @{ Row = 13; B = 107.83000199999999; BFmt = $true; C = 86.900002000000001; CFmt = $true; D = 0.27; E = 363 },
    @{ Row = 14; B = 140.52000000000001; BFmt = $false; C = 102.42; CFmt = $false; D = 0.27; E = 622 },
    @{ Row = 15; B = 336.77; BFmt = $false; C = 223.92; CFmt = $false; D = 0.33500000000000002; E = 101 },
    @{ Row = 16; B = 368.95; BFmt = $false; C = 295.45999999999998; CFmt = $false; D = 0.19900000000000001; E = 87 },
    @{ Row = 17; B = 1186.75; BFmt = $false; C = 957.28; CFmt = $false; D = 0.193; E = 45 },
    @{ Row = 18; B = 1527.46; BFmt = $false; C = 776.76; CFmt = $false; D = 0.49099999999999999; E = 929 },
    @{ Row = 19; B = 1565.15; BFmt = $false; C = 676.53; CFmt = $false; D = 0.56799999999999995; E = 517 },
    @{ Row = 20; B = 1360.48; BFmt = $false; C = 1099.23; CFmt = $false; D = 0.21579999999999999; E = 155 },
    @{ Row = 21; B = 2930.75; BFmt = $false; C = 2351.1; CFmt = $false; D = 0.19800000000000001; E = 95 },
    @{ Row = 22; B = 3386.15; BFmt = $false; C = 2237.4; CFmt = $false; D = 0.33900000000000002; E = 33 },
    @{ Row = 23; B = 4796.5600000000004; BFmt = $false; C = 3577.03; CFmt = $false; D = 0.27550000000000002; E = 282 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    if ($r.BFmt) { $ws.Cells.Item($r.Row, 2).NumberFormat = "0.00" }
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    if ($r.CFmt) { $ws.Cells.Item($r.Row, 3).NumberFormat = "0.00" }
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 4).NumberFormat = "0.00%"
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# --- Step 3: sheet selection.
$ws.Range("H4").Select()

# --- Step 4: workbook window size/position.
$excel.Width = 33040
$excel.Height = 15580
